# Update odds values on Sheet1 to reflect the latest FlashScore scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Aston Villa vs Liverpool)
$ws.Range("O2").Value = 1.2
$ws.Range("P2").Value = 4.5
$ws.Range("Q2").Value = 1.62
$ws.Range("R2").Value = 2.3
$ws.Range("S2").Value = 2.01
$ws.Range("T2").Value = 1.89

# Row 3 (Once Caldas vs Pereira)
$ws.Range("G3").Value = 1.95
$ws.Range("H3").Value = 3.25
$ws.Range("J3").Value = 2.75
$ws.Range("K3").Value = 1.95
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("Q3").Value = 2.4
$ws.Range("R3").Value = 1.53
$ws.Range("U3").Value = 5
$ws.Range("V3").Value = 1.17
$ws.Range("Y3").Value = 2.2
$ws.Range("Z3").Value = 1.62
$ws.Range("AC3").Value = 9.5
$ws.Range("AD3").Value = 17
$ws.Range("AF3").Value = 41
$ws.Range("AI3").Value = 21
$ws.Range("AJ3").Value = 81
$ws.Range("AL3").Value = 9
$ws.Range("AM3").Value = 19
$ws.Range("AR3").Value = 1.88
$ws.Range("AS3").Value = 1.98

# Row 4 (Luton vs Plymouth)
$ws.Range("G4").Value = 1.75
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 2.4
$ws.Range("K4").Value = 2.05
$ws.Range("U4").Value = 4
$ws.Range("V4").Value = 1.25
$ws.Range("AB4").Value = 7.5
$ws.Range("AC4").Value = 9
$ws.Range("AE4").Value = 17
$ws.Range("AH4").Value = 6.5

# Row 5 (Hyderabad vs Mumbai City)
$ws.Range("H5").Value = 3.9
$ws.Range("I5").Value = 1.62
$ws.Range("K5").Value = 2.38
$ws.Range("O5").Value = 1.2
$ws.Range("P5").Value = 4.33
$ws.Range("Q5").Value = 1.65
$ws.Range("R5").Value = 2.2
$ws.Range("S5").Value = 2.03
$ws.Range("T5").Value = 1.78
$ws.Range("U5").Value = 2.63
$ws.Range("V5").Value = 1.44
$ws.Range("AF5").Value = 41
$ws.Range("AL5").Value = 8.5
$ws.Range("AM5").Value = 8.5

# Row 7 (Club America vs Club Leon)
$ws.Range("G7").Value = 1.56
$ws.Range("AB7").Value = 8
$ws.Range("AE7").Value = 12
